$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2: target cluster becomes "ECs", recompute derived metrics ---
$ws.Cells.Item(2, 4).Value = "ECs"

$ws.Cells.Item(2, 13).Value = 0.014278
$ws.Cells.Item(2, 14).Value = 0.042834
$ws.Cells.Item(2, 15).Value = 0.08545683615785395
$ws.Cells.Item(2, 16).Value = 0.109974299791266
$ws.Cells.Item(2, 17).Value = 0.01644934588733333
$ws.Cells.Item(2, 18).Value = 0.148044112986
$ws.Cells.Item(2, 19).Value = 0.08545683615785395
$ws.Cells.Item(2, 20).Value = 0.109974299791266

# --- New row 3: FAPs -> FAPs (original values that used to live on row 2) ---
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Rln3"
$ws.Cells.Item(3, 3).Value = "Rxfp1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.152076333333333
$ws.Cells.Item(3, 8).Value = 3.456229
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.041056
$ws.Cells.Item(3, 14).Value = 0.123168
$ws.Cells.Item(3, 15).Value = 0.2457288041250071
$ws.Cells.Item(3, 16).Value = 0.3162281028316444
$ws.Cells.Item(3, 17).Value = 0.04729964594133333
$ws.Cells.Item(3, 18).Value = 0.425696813472
$ws.Cells.Item(3, 19).Value = 0.2457288041250071
$ws.Cells.Item(3, 20).Value = 0.3162281028316444

# --- New row 4: FAPs -> MuSCs ---
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Rln3"
$ws.Cells.Item(4, 3).Value = "Rxfp1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.152076333333333
$ws.Cells.Item(4, 8).Value = 3.456229
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.1117445
$ws.Cells.Item(4, 14).Value = 0.223489
$ws.Cells.Item(4, 15).Value = 0.668814359717139
$ws.Cells.Item(4, 16).Value = 0.5737975973770896
$ws.Cells.Item(4, 17).Value = 0.1287381938301667
$ws.Cells.Item(4, 18).Value = 0.772429162981
$ws.Cells.Item(4, 19).Value = 0.668814359717139
$ws.Cells.Item(4, 20).Value = 0.5737975973770896
